# Update FHIR StructureDefinition spreadsheet from Alvearie/IBM -> LinuxForHealth
# branding, bump the version/date, and refresh the example ValueSet URLs.

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/condition-disease-course"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Extension.url fixed value now points at the renamed StructureDefinition URL
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/condition-disease-course"

# Extension.value[x] (valueCodeableConcept slice) example binding value set
$elements.Range("Y7").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/condition-course"

# The root Extension row no longer carries the ele-1/ext-1 constraint text
$elements.Range("AI2").Value = ""

# Column Y widens to fit the longer Binding Value Set URL
$elements.Columns.Item(25).ColumnWidth = 55.98828125
